# Add a new worksheet "13.07, 41 deg, 48 hr" at the end of the workbook.
#
# The new sheet's well/treatment layout is identical to the "11.07, 41 deg"
# sheet (same well labels in column A, same fRS585 wells in rows 2-9), except
# the "blank" wells (rows 10-17) are labelled "Blank" here instead of the
# plain-text "blank" used on the source sheet.
#
# We build it by copying "11.07, 41 deg" (so we inherit its number format,
# cell style and drawing placeholder), then:
#   - drop the leftover column D formatting-only cells that belonged to the
#     source sheet's comment anchor,
#   - remove the copied cell comment (and its legacy VML drawing),
#   - update column B for rows 10-17 from "blank" to "Blank",
#   - rename the sheet and make sure it is the last tab.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("11.07, 41 deg")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy the source sheet so we inherit formatting/drawing, placing the copy
# right after the current last sheet.
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "13.07, 41 deg, 48 hr"

# The source sheet has a couple of style-only cells in column D (anchors for
# its comment's formatting); the new sheet doesn't have any of that.
$newSheet.Range("D1:D17").Clear()

# Drop the copied comment (and its legacy VML drawing) - the new sheet has
# none.
$commentCell = $newSheet.Range("A1").Comment
if ($commentCell -ne $null) {
    $commentCell.Delete()
}

# Rows 10-17 in column B: "blank" -> "Blank"
$newSheet.Range("B10:B17").Value = "Blank"

# Make sure the new sheet ends up as the very last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
if ($lastSheet.Name -ne "13.07, 41 deg, 48 hr") {
    $newSheet.Move($null, $lastSheet)
}
